# Update the "想去人数" (interested-count) figures to reflect newly scraped
# totals for this gh-pages data refresh (commit 456a3b4).
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition listing)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F8").Value  = 3651   # 上海·趣元界-INW动漫游戏展
$ws1.Range("F9").Value  = 1257   # 上海·Nova次元动漫嘉年华
$ws1.Range("F27").Value = 406    # 上海·夜蓝诗2.0·恋与深空同人only
$ws1.Range("F30").Value = 86     # 上海·第五届长三角文博会上海国际插画艺术节
$ws1.Range("F31").Value = 86     # 上海·第五届长三角文博会上海国际插画艺术节
$ws1.Range("F33").Value = 2690   # 上海·第二届iPR动漫游戏嘉年华（取消）
$ws1.Range("F34").Value = 157    # 上海·明日方舟同人ONLY
$ws1.Range("F36").Value = 1229   # 上海·向前冲！运动番同人Only

# Sheet "本地生活" (Local-life listing)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F11").Value = 567    # 上海·三丽鸥家族Sanrio Characters主题餐厅·海滩奇遇季

# Sheet "全部类型" (Combined/all-types listing)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value  = 567    # 上海·三丽鸥家族Sanrio Characters主题餐厅·海滩奇遇季
$ws4.Range("F12").Value = 1257   # 上海·Nova次元动漫嘉年华
$ws4.Range("F34").Value = 406    # 上海·夜蓝诗2.0·恋与深空同人only
$ws4.Range("F40").Value = 86     # 上海·第五届长三角文博会上海国际插画艺术节
$ws4.Range("F45").Value = 157    # 上海·明日方舟同人ONLY
$ws4.Range("F49").Value = 1229   # 上海·向前冲！运动番同人Only
